# Applies the "Updated cryptos list" crypto price/volume refresh to sheet1.
# For Price (column D) values that are purely numeric-looking (e.g. "568.87"),
# a leading apostrophe is used so Excel stores/keeps them as text, matching
# the original inlineStr cell content instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.895.58"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "2.897.66"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'568.87"
$ws.Range("E5").Value = "  -3.28%  "

$ws.Range("D6").Value = "'143.86"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'0.501"
$ws.Range("E8").Value = "  -0.93%  "

$ws.Range("D9").Value = "2.895.54"
$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("D10").Value = "'6.89"
$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").Value = "'0.0000231"
$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("D14").Value = "'32.37"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "'0.126"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").Value = "3.378.75"
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("D17").Value = "61.848.72"
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.897.70"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.52"
$ws.Range("E19").Value = "  -1.78%  "

$ws.Range("D20").Value = "'431.37"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").Value = "'12.94"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("D24").Value = "'78.83"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("D26").Value = "'10.09"
$ws.Range("E26").Value = "  -8.72%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D29").Value = "'0.0000111"
$ws.Range("E29").Value = "  +10.96%  "

$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("D31").Value = "'2.51"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("E32").Value = "  -5.06%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").Value = "'25.53"
$ws.Range("E35").Value = "  -2.28%  "

$ws.Range("E36").Value = "  -3.71%  "

$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("E39").Value = "  -4.73%  "

$ws.Range("E40").Value = "  -4.16%  "

$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("E42").Value = "  -2.26%  "

$ws.Range("D43").Value = "'40.25"
$ws.Range("E43").Value = "  +4.28%  "

$ws.Range("D44").Value = "'0.269"
$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").Value = "2.697.86"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").Value = "'131.75"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("D48").Value = "'347.69"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "'21.62"
$ws.Range("E51").Value = "  -3.47%  "
